# Scheduled runner update: refresh market-board derived profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit* columns) for a batch of
# leve rows across the crafting job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(69, 8).Value = 3576.375
$ws.Cells.Item(69, 9).Value = 3362.4
$ws.Cells.Item(69, 10).Value = 3933
$ws.Cells.Item(69, 11).Value = 10087.2
$ws.Cells.Item(69, 12).Value = 11799
$ws.Cells.Item(69, 13).Value = -9213.200000000001
$ws.Cells.Item(69, 14).Value = -13547

$ws.Cells.Item(72, 8).Value = 3576.375
$ws.Cells.Item(72, 9).Value = 3362.4
$ws.Cells.Item(72, 10).Value = 3933
$ws.Cells.Item(72, 11).Value = 30261.6
$ws.Cells.Item(72, 12).Value = 35397
$ws.Cells.Item(72, 13).Value = -25893.6
$ws.Cells.Item(72, 14).Value = -44133

$ws.Cells.Item(124, 8).Value = 50592.668
$ws.Cells.Item(124, 10).Value = 50592.668
$ws.Cells.Item(124, 12).Value = 50592.668
$ws.Cells.Item(124, 14).Value = -60412.668

$ws.Cells.Item(128, 8).Value = 46772
$ws.Cells.Item(128, 10).Value = 46772
$ws.Cells.Item(128, 12).Value = 46772
$ws.Cells.Item(128, 14).Value = -56732

$ws.Cells.Item(130, 8).Value = 37193.6
$ws.Cells.Item(130, 10).Value = 43992
$ws.Cells.Item(130, 12).Value = 43992
$ws.Cells.Item(130, 14).Value = -54032

$ws.Cells.Item(132, 8).Value = 24572.738
$ws.Cells.Item(132, 9).Value = 3862.182
$ws.Cells.Item(132, 10).Value = 100511.445
$ws.Cells.Item(132, 11).Value = 11586.546
$ws.Cells.Item(132, 12).Value = 301534.335
$ws.Cells.Item(132, 13).Value = -9056.545999999998
$ws.Cells.Item(132, 14).Value = -306594.335

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(31, 8).Value = 11421.5
$ws.Cells.Item(31, 9).Value = 2843
$ws.Cells.Item(31, 11).Value = 2843
$ws.Cells.Item(31, 13).Value = -2549

$ws.Cells.Item(32, 8).Value = 42791.445
$ws.Cells.Item(32, 9).Value = 42553.87
$ws.Cells.Item(32, 10).Value = 44335.668
$ws.Cells.Item(32, 11).Value = 42553.87
$ws.Cells.Item(32, 12).Value = 44335.668
$ws.Cells.Item(32, 13).Value = -42266.87
$ws.Cells.Item(32, 14).Value = -44909.668

$ws.Cells.Item(61, 8).Value = 3248.2273
$ws.Cells.Item(61, 9).Value = 1990.381
$ws.Cells.Item(61, 10).Value = 4396.696
$ws.Cells.Item(61, 11).Value = 1990.381
$ws.Cells.Item(61, 12).Value = 4396.696
$ws.Cells.Item(61, 13).Value = -1778.381
$ws.Cells.Item(61, 14).Value = -4820.696

$ws.Cells.Item(74, 8).Value = 1520.25
$ws.Cells.Item(74, 9).Value = 879.13336
$ws.Cells.Item(74, 10).Value = 2588.7778
$ws.Cells.Item(74, 11).Value = 879.13336
$ws.Cells.Item(74, 12).Value = 2588.7778
$ws.Cells.Item(74, 13).Value = -5.133360000000039
$ws.Cells.Item(74, 14).Value = -4336.7778

$ws.Cells.Item(77, 8).Value = 1520.25
$ws.Cells.Item(77, 9).Value = 879.13336
$ws.Cells.Item(77, 10).Value = 2588.7778
$ws.Cells.Item(77, 11).Value = 4395.6668
$ws.Cells.Item(77, 12).Value = 12943.889
$ws.Cells.Item(77, 13).Value = -27.66679999999997
$ws.Cells.Item(77, 14).Value = -21679.889

$ws.Cells.Item(123, 8).Value = 51429
$ws.Cells.Item(123, 10).Value = 51429
$ws.Cells.Item(123, 12).Value = 51429
$ws.Cells.Item(123, 14).Value = -61229

$ws.Cells.Item(136, 8).Value = 3248.2273
$ws.Cells.Item(136, 9).Value = 1990.381
$ws.Cells.Item(136, 10).Value = 4396.696
$ws.Cells.Item(136, 11).Value = 5971.143
$ws.Cells.Item(136, 12).Value = 13190.088
$ws.Cells.Item(136, 13).Value = -3421.143
$ws.Cells.Item(136, 14).Value = -18290.088

$ws.Cells.Item(138, 8).Value = 42500
$ws.Cells.Item(138, 10).Value = 42500
$ws.Cells.Item(138, 12).Value = 42500
$ws.Cells.Item(138, 14).Value = -52780

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(139, 8).Value = 55166.668
$ws.Cells.Item(139, 10).Value = 55166.668
$ws.Cells.Item(139, 12).Value = 55166.668
$ws.Cells.Item(139, 14).Value = -65446.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(51, 8).Value = 125029840
$ws.Cells.Item(51, 10).Value = 39786.332
$ws.Cells.Item(51, 12).Value = 39786.332
$ws.Cells.Item(51, 14).Value = -41258.332

$ws.Cells.Item(58, 8).Value = 1930.4062
$ws.Cells.Item(58, 9).Value = 1759.2174
$ws.Cells.Item(58, 10).Value = 2367.889
$ws.Cells.Item(58, 11).Value = 1759.2174
$ws.Cells.Item(58, 12).Value = 2367.889
$ws.Cells.Item(58, 13).Value = -1556.2174
$ws.Cells.Item(58, 14).Value = -2773.889

$ws.Cells.Item(61, 8).Value = 125029840
$ws.Cells.Item(61, 10).Value = 39786.332
$ws.Cells.Item(61, 12).Value = 39786.332
$ws.Cells.Item(61, 14).Value = -40482.332

$ws.Cells.Item(136, 8).Value = 1930.4062
$ws.Cells.Item(136, 9).Value = 1759.2174
$ws.Cells.Item(136, 10).Value = 2367.889
$ws.Cells.Item(136, 11).Value = 5277.6522
$ws.Cells.Item(136, 12).Value = 7103.667
$ws.Cells.Item(136, 13).Value = -2727.6522
$ws.Cells.Item(136, 14).Value = -12203.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(86, 8).Value = 1708.8182
$ws.Cells.Item(86, 9).Value = 800
$ws.Cells.Item(86, 10).Value = 1910.7778
$ws.Cells.Item(86, 11).Value = 2400
$ws.Cells.Item(86, 12).Value = 5732.3334
$ws.Cells.Item(86, 13).Value = -1214
$ws.Cells.Item(86, 14).Value = -8104.3334

$ws.Cells.Item(89, 8).Value = 1708.8182
$ws.Cells.Item(89, 9).Value = 800
$ws.Cells.Item(89, 10).Value = 1910.7778
$ws.Cells.Item(89, 11).Value = 7200
$ws.Cells.Item(89, 12).Value = 17197.0002
$ws.Cells.Item(89, 13).Value = -1272
$ws.Cells.Item(89, 14).Value = -29053.0002

$ws.Cells.Item(113, 8).Value = 4502.385
$ws.Cells.Item(113, 9).Value = 6902.75
$ws.Cells.Item(113, 10).Value = 661.8
$ws.Cells.Item(113, 11).Value = 20708.25
$ws.Cells.Item(113, 12).Value = 1985.4
$ws.Cells.Item(113, 13).Value = -18538.25
$ws.Cells.Item(113, 14).Value = -6325.4

$ws.Cells.Item(117, 8).Value = 1032.25

$ws.Cells.Item(136, 8).Value = 29414226
$ws.Cells.Item(136, 9).Value = 38463444
$ws.Cells.Item(136, 10).Value = 4266.5
$ws.Cells.Item(136, 11).Value = 115390332
$ws.Cells.Item(136, 12).Value = 12799.5
$ws.Cells.Item(136, 13).Value = -115385232
$ws.Cells.Item(136, 14).Value = -22999.5

$ws.Cells.Item(138, 8).Value = 3617.6553
$ws.Cells.Item(138, 9).Value = 3145.84
$ws.Cells.Item(138, 10).Value = 6566.5
$ws.Cells.Item(138, 11).Value = 9437.52
$ws.Cells.Item(138, 12).Value = 19699.5
$ws.Cells.Item(138, 13).Value = -4297.52
$ws.Cells.Item(138, 14).Value = -29979.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 176026.97
$ws.Cells.Item(80, 9).Value = 316311.25
$ws.Cells.Item(80, 10).Value = 3369.3845
$ws.Cells.Item(80, 11).Value = 316311.25
$ws.Cells.Item(80, 12).Value = 3369.3845
$ws.Cells.Item(80, 13).Value = -315313.25
$ws.Cells.Item(80, 14).Value = -5365.3845

$ws.Cells.Item(83, 8).Value = 176026.97
$ws.Cells.Item(83, 9).Value = 316311.25
$ws.Cells.Item(83, 10).Value = 3369.3845
$ws.Cells.Item(83, 11).Value = 1581556.25
$ws.Cells.Item(83, 12).Value = 16846.9225
$ws.Cells.Item(83, 13).Value = -1576564.25
$ws.Cells.Item(83, 14).Value = -26830.9225

$ws.Cells.Item(130, 8).Value = 52964.8
$ws.Cells.Item(130, 10).Value = 52964.8
$ws.Cells.Item(130, 12).Value = 52964.8
$ws.Cells.Item(130, 14).Value = -63004.8

$ws.Cells.Item(138, 8).Value = 48000
$ws.Cells.Item(138, 10).Value = 48000
$ws.Cells.Item(138, 12).Value = 48000
$ws.Cells.Item(138, 14).Value = -58280

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(137, 8).Value = 50996.668
$ws.Cells.Item(137, 10).Value = 50996.668
$ws.Cells.Item(137, 12).Value = 50996.668
$ws.Cells.Item(137, 14).Value = -61196.668

$ws.Cells.Item(138, 8).Value = 44077.5
$ws.Cells.Item(138, 10).Value = 44077.5
$ws.Cells.Item(138, 12).Value = 44077.5
$ws.Cells.Item(138, 14).Value = -54357.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(119, 8).Value = 43857.6
$ws.Cells.Item(119, 10).Value = 43857.6
$ws.Cells.Item(119, 12).Value = 43857.6
$ws.Cells.Item(119, 14).Value = -53533.6

$ws.Cells.Item(137, 8).Value = 47178.5
$ws.Cells.Item(137, 10).Value = 47178.5
$ws.Cells.Item(137, 12).Value = 47178.5
$ws.Cells.Item(137, 14).Value = -57378.5

$ws.Cells.Item(139, 8).Value = 28113
$ws.Cells.Item(139, 10).Value = 28113
$ws.Cells.Item(139, 12).Value = 28113
$ws.Cells.Item(139, 14).Value = -38393
